# Update the cryptocurrency price/volume table (cols D and E) with the
# latest scraped values, as produced by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.413.44"
$ws.Range("E2").Value = "  +4.55%  "
$ws.Range("D3").Value = "4.061.41"
$ws.Range("E3").Value = "  +4.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.13"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.20"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  +18.91%  "
$ws.Range("D8").Value = "4.054.19"
$ws.Range("E8").Value = "  +4.45%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +7.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.176"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000329"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.16"
$ws.Range("E13").Value = "  +14.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.21"
$ws.Range("E14").Value = "  +9.15%  "
$ws.Range("D15").Value = "4.705.47"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D16").Value = "4.076.57"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.38"
$ws.Range("E17").Value = "  +8.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.20"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").Value = "72.376.66"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.52"
$ws.Range("E22").Value = "  +4.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "104.53"
$ws.Range("E23").Value = "  +18.42%  "
$ws.Range("E24").Value = "  +6.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.92"
$ws.Range("E25").Value = "  +5.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.02"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.64"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.08"
$ws.Range("E28").Value = "  +4.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.87"
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("E31").Value = "  +15.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.73"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("E33").Value = "  +4.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "681.92"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.87"
$ws.Range("E35").Value = "  +14.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.15"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "43.35"
$ws.Range("E37").Value = "  +8.43%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "0.0₃0864"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  +9.53%  "
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0499"
$ws.Range("E43").Value = "  +3.87%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.157"
$ws.Range("E46").Value = "  +12.16%  "
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").Value = "  +3.93%  "
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("E50").Value = "  +7.00%  "
$ws.Range("E51").Value = "  +2.42%  "
